# [EXTRA SCRAPE] full data scraped for extra batting and bowling fields
#
# 1) "ODI Batting" sheet: several rows had an empty INNING_NUMBER (column B)
#    cell that should not exist at all any more (the inning number for those
#    matches is simply unknown/not applicable). Remove the empty cells.
# 2) Add a brand-new "ODI Batting Extra" sheet (after "ODI Bowling") holding
#    additional scraped batting fields per match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Clear the stray empty INNING_NUMBER cells on the "ODI Batting" sheet
# ---------------------------------------------------------------------------
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$emptyInningRows = @(2, 3, 8, 9, 11, 12, 16)
foreach ($r in $emptyInningRows) {
    $odiBatting.Cells.Item($r, 2).ClearContents()
}

# ---------------------------------------------------------------------------
# 2. Add the new "ODI Batting Extra" sheet at the end of the workbook
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ODI Batting Extra"

# Match the page setup used by the rest of the workbook's sheets.
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36
$ws.Outline.SummaryBelow = $true
$ws.Outline.SummaryRight = $true

# Header row - reuse the exact header style already used by the other sheets
# (bold font, thin border, centered) by copying its formatting across.
$odiBatting.Range("A1:D1").Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)

$ws.Cells.Item(1, 1).Value = "MATCH_CODE"
$ws.Cells.Item(1, 2).Value = "BATTING_POSITION"
$ws.Cells.Item(1, 3).Value = "NUM_4"
$ws.Cells.Item(1, 4).Value = "NUM_6"
$ws.Cells.Item(1, 5).Value = "PERCENT_RUNS_OF_TOTAL"
$ws.Cells.Item(1, 6).Value = "MAN_OF_MATCH"

# Columns A, C, D, E hold scraped text values (match codes, counts and a
# percentage string) that must stay as literal text rather than being
# auto-converted to numbers/percentages by Excel's type inference.
$ws.Range("A2:A22").NumberFormat = "@"
$ws.Range("C2:E22").NumberFormat = "@"

# Data rows 2-22
$ws.Cells.Item(2,1).Value = '3340'
$ws.Cells.Item(2,2).Value = 9
$ws.Cells.Item(2,6).Value = 'NO'

$ws.Cells.Item(3,1).Value = '3342'
$ws.Cells.Item(3,6).Value = 'NO'

$ws.Cells.Item(4,1).Value = '3366'
$ws.Cells.Item(4,2).Value = 9
$ws.Cells.Item(4,3).Value = '0'
$ws.Cells.Item(4,4).Value = '0'
$ws.Cells.Item(4,5).Value = '3.23%'
$ws.Cells.Item(4,6).Value = 'NO'

$ws.Cells.Item(5,1).Value = '3370'
$ws.Cells.Item(5,6).Value = 'NO'

$ws.Cells.Item(6,1).Value = '3385'
$ws.Cells.Item(6,6).Value = 'NO'

$ws.Cells.Item(7,1).Value = '3429'
$ws.Cells.Item(7,2).Value = 9
$ws.Cells.Item(7,3).Value = '0'
$ws.Cells.Item(7,4).Value = '0'
$ws.Cells.Item(7,5).Value = '2.08%'
$ws.Cells.Item(7,6).Value = 'NO'

$ws.Cells.Item(8,1).Value = '3503'
$ws.Cells.Item(8,2).Value = 9
$ws.Cells.Item(8,6).Value = 'NO'

$ws.Cells.Item(9,1).Value = '3827'
$ws.Cells.Item(9,2).Value = 8
$ws.Cells.Item(9,6).Value = 'NO'

$ws.Cells.Item(10,1).Value = '3828'
$ws.Cells.Item(10,2).Value = 8
$ws.Cells.Item(10,3).Value = '1'
$ws.Cells.Item(10,4).Value = '0'
$ws.Cells.Item(10,5).Value = '5.88%'
$ws.Cells.Item(10,6).Value = 'NO'

$ws.Cells.Item(11,1).Value = '3865'
$ws.Cells.Item(11,2).Value = 8
$ws.Cells.Item(11,6).Value = 'NO'

$ws.Cells.Item(12,1).Value = '3866'
$ws.Cells.Item(12,2).Value = 8
$ws.Cells.Item(12,6).Value = 'NO'

$ws.Cells.Item(13,1).Value = '3868'
$ws.Cells.Item(13,6).Value = 'NO'

$ws.Cells.Item(14,1).Value = '3888'
$ws.Cells.Item(14,2).Value = 8
$ws.Cells.Item(14,3).Value = '0'
$ws.Cells.Item(14,4).Value = '0'
$ws.Cells.Item(14,5).Value = '0.81%'
$ws.Cells.Item(14,6).Value = 'NO'

$ws.Cells.Item(15,1).Value = '3951'
$ws.Cells.Item(15,6).Value = 'NO'

$ws.Cells.Item(16,1).Value = '4100'
$ws.Cells.Item(16,2).Value = 8
$ws.Cells.Item(16,6).Value = 'YES'

$ws.Cells.Item(17,1).Value = '4101'
$ws.Cells.Item(17,2).Value = 8
$ws.Cells.Item(17,3).Value = '1'
$ws.Cells.Item(17,4).Value = '0'
$ws.Cells.Item(17,5).Value = '1.54%'
$ws.Cells.Item(17,6).Value = 'NO'

$ws.Cells.Item(18,1).Value = '4239'
$ws.Cells.Item(18,2).Value = 8
$ws.Cells.Item(18,3).Value = '1'
$ws.Cells.Item(18,4).Value = '0'
$ws.Cells.Item(18,5).Value = '4.46%'
$ws.Cells.Item(18,6).Value = 'NO'

$ws.Cells.Item(19,1).Value = '4242'
$ws.Cells.Item(19,6).Value = 'NO'

$ws.Cells.Item(20,1).Value = '4245'
$ws.Cells.Item(20,2).Value = 8
$ws.Cells.Item(20,3).Value = '1'
$ws.Cells.Item(20,4).Value = '0'
$ws.Cells.Item(20,5).Value = '6.17%'
$ws.Cells.Item(20,6).Value = 'NO'

$ws.Cells.Item(21,1).Value = '4566'
$ws.Cells.Item(21,2).Value = 8
$ws.Cells.Item(21,3).Value = '2'
$ws.Cells.Item(21,4).Value = '0'
$ws.Cells.Item(21,5).Value = '15.53%'
$ws.Cells.Item(21,6).Value = 'NO'

$ws.Cells.Item(22,1).Value = '4568'
$ws.Cells.Item(22,6).Value = 'NO'
